$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Enho"
$ws.Range("C2").Value = "Gpr19"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.193526
$ws.Range("H2").Value = 0.580578
$ws.Range("I2").Value = 0.1547204586252329
$ws.Range("J2").Value = 0.1547204586252329
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.7909473333333334
$ws.Range("N2").Value = 2.372842
$ws.Range("O2").Value = 0.1281079032992492
$ws.Range("P2").Value = 0.1281079032992493
$ws.Range("Q2").Value = 0.1530688736306667
$ws.Range("R2").Value = 1.377619862676
$ws.Range("S2").Value = 0.01982091355197683
$ws.Range("T2").Value = 0.01982091355197683

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Enho"
$ws.Range("C3").Value = "Gpr19"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.193526
$ws.Range("H3").Value = 0.580578
$ws.Range("I3").Value = 0.1547204586252329
$ws.Range("J3").Value = 0.1547204586252329
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.366842
$ws.Range("N3").Value = 7.100526
$ws.Range("O3").Value = 0.3833519038274798
$ws.Range("P3").Value = 0.3833519038274799
$ws.Range("Q3").Value = 0.458045464892
$ws.Range("R3").Value = 4.122409184028
$ws.Range("S3").Value = 0.05931238237504385
$ws.Range("T3").Value = 0.05931238237504387

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Enho"
$ws.Range("C4").Value = "Gpr19"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.193526
$ws.Range("H4").Value = 0.580578
$ws.Range("I4").Value = 0.1547204586252329
$ws.Range("J4").Value = 0.1547204586252329
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.016282
$ws.Range("N4").Value = 9.048845999999999
$ws.Range("O4").Value = 0.4885401928732709
$ws.Range("P4").Value = 0.488540192873271
$ws.Range("Q4").Value = 0.583728990332
$ws.Range("R4").Value = 5.253560912988
$ws.Range("S4").Value = 0.07558716269821221
$ws.Range("T4").Value = 0.07558716269821224

$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Enho"
$ws.Range("C5").Value = "Gpr19"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.057284666666667
$ws.Range("H5").Value = 3.171854
$ws.Range("I5").Value = 0.845279541374767
$ws.Range("J5").Value = 0.845279541374767
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7909473333333334
$ws.Range("N5").Value = 2.372842
$ws.Range("O5").Value = 0.1281079032992492
$ws.Range("P5").Value = 0.1281079032992493
$ws.Range("Q5").Value = 0.8362564876742223
$ws.Range("R5").Value = 7.526308389068
$ws.Range("S5").Value = 0.1082869897472724
$ws.Range("T5").Value = 0.1082869897472724

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Enho"
$ws.Range("C6").Value = "Gpr19"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.057284666666667
$ws.Range("H6").Value = 3.171854
$ws.Range("I6").Value = 0.845279541374767
$ws.Range("J6").Value = 0.845279541374767
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.366842
$ws.Range("N6").Value = 7.100526
$ws.Range("O6").Value = 0.3833519038274798
$ws.Range("P6").Value = 0.3833519038274799
$ws.Range("Q6").Value = 2.502425755022667
$ws.Range("R6").Value = 22.521831795204
$ws.Range("S6").Value = 0.3240395214524359
$ws.Range("T6").Value = 0.324039521452436

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Enho"
$ws.Range("C7").Value = "Gpr19"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.057284666666667
$ws.Range("H7").Value = 3.171854
$ws.Range("I7").Value = 0.845279541374767
$ws.Range("J7").Value = 0.845279541374767
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.016282
$ws.Range("N7").Value = 9.048845999999999
$ws.Range("O7").Value = 0.4885401928732709
$ws.Range("P7").Value = 0.488540192873271
$ws.Range("Q7").Value = 3.189068708942667
$ws.Range("R7").Value = 28.701618380484
$ws.Range("S7").Value = 0.4129530301750587
$ws.Range("T7").Value = 0.4129530301750587

# Remove now-obsolete rows 8-10 (ECs as sending cluster rows removed)
$ws.Range("A8:T10").EntireRow.Delete()
